$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Ativação:" date value 01/01/2012 -> 01/01/2023 ---
# Use a scratch cell + PasteSpecial(values) so the text is stored as a
# literal string (shared string) instead of being auto-parsed into a date
# serial number, and PasteSpecial(formats) below preserves original styles
# for brand-new cells.
$ws.Range("E1").Formula = "=""01/01/2023"""
$ws.Range("E1").Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Range("C8").PasteSpecial(-4163)
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("E1").Clear()

# --- Row 11 (Objectives:) - add English objectives text in B/C ---
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("B11").Value = "Provide knowledge about vacuum systems and production techniques and use of low temperatures."
$ws.Range("C11").Value = "Provide knowledge about vacuum systems and production techniques and use of low temperatures."

# --- Row 14 (Short syllabus:) - add short syllabus text in B/C ---
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("B14").Value = "Vacuum systems. Cryogenics and low temperature."
$ws.Range("C14").Value = "Vacuum systems. Cryogenics and low temperature."

# --- Row 16 (Syllabus:) - add full syllabus text in B/C ---
$ws.Range("B13").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("B16").Value = "Theory of rarefied gases. Gas flow. Vacuum pumps. Quantitative description of the pumping of vacuum systems. Pressure gauges. Accessories: traps, shields, valves, etc. Adsorption, desorption and evaporation of molecules in vacuum. Leak detection .Sealing.Welding.Cleaning.cryogenics. Properties of cryogenic gases and liquids. Methods for obtaining low temperature. Liquefaction of gases. Temperature measurement. Cryogenic components. Calculation of heat transfer in cryostats and dewars."
$ws.Range("C16").Value = "Theory of rarefied gases. Gas flow. Vacuum pumps. Quantitative description of the pumping of vacuum systems. Pressure gauges. Accessories: traps, shields, valves, etc. Adsorption, desorption and evaporation of molecules in vacuum. Leak detection .Sealing.Welding.Cleaning.cryogenics. Properties of cryogenic gases and liquids. Methods for obtaining low temperature. Liquefaction of gases. Temperature measurement. Cryogenic components. Calculation of heat transfer in cryostats and dewars."
